# Update the date line and the twenty-five multiplication problems
# to match the new day's generated content.

$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-27 Tuesday", "2026-01-28 Wednesday"),
    @("59×93=", "63×22="),
    @("95×27=", "71×37="),
    @("36×21=", "91×15="),
    @("51×35=", "17×86="),
    @("87×97=", "73×75="),
    @("16×67=", "72×68="),
    @("22×97=", "85×33="),
    @("22×70=", "45×17="),
    @("89×14=", "80×13="),
    @("58×42=", "78×97="),
    @("46×29=", "79×84="),
    @("20×63=", "62×69="),
    @("52×97=", "59×26="),
    @("16×46=", "44×27="),
    @("20×12=", "86×45="),
    @("76×48=", "66×14="),
    @("36×78=", "62×84="),
    @("56×40=", "64×31="),
    @("37×38=", "12×72="),
    @("53×18=", "31×33="),
    @("78×25=", "85×83="),
    @("35×39=", "57×28="),
    @("13×71=", "53×88="),
    @("37×58=", "95×13="),
    @("71×72=", "87×60=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
    Write-Host "Replaced '$old' -> '$new'"
}
